$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = 5912805
$ws.Range("C13").Value = 159110
$ws.Range("B14").Value = 1807236
$ws.Range("C14").Value = 9642307
$ws.Range("B15").Value = 7449809
$ws.Range("C15").Value = 2294202
$ws.Range("B16").Value = 9206617
$ws.Range("C16").Value = 7057512
